$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8110.7144
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 8962.5
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 26887.5
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -28635.5
$ws.Range("H72").Value = 8110.7144
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 8962.5
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 80662.5
$ws.Range("M72").Value = -22632
$ws.Range("N72").Value = -89398.5
$ws.Range("H76").Value = 3114.2856
$ws.Range("I76").Value = 3114.2856
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3114.2856
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -2799.2856
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 3114.2856
$ws.Range("I79").Value = 3114.2856
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3114.2856
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2022.2856
$ws.Range("N79").ClearContents()
$ws.Range("H118").Value = 545
$ws.Range("I118").Value = 315
$ws.Range("J118").Value = 775
$ws.Range("K118").Value = 945
$ws.Range("L118").Value = 2325
$ws.Range("M118").Value = 712
$ws.Range("N118").Value = -5639
$ws.Range("H127").Value = 111111660
$ws.Range("I127").Value = 111111660
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 333334980
$ws.Range("L127").Value = 0
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -333330020
$ws.Range("H129").Value = 894.1
$ws.Range("I129").Value = 270.66666
$ws.Range("J129").Value = 1049.9584
$ws.Range("K129").Value = 811.9999799999999
$ws.Range("L129").Value = 3149.8752
$ws.Range("M129").Value = 4188.00002
$ws.Range("N129").Value = -13149.8752
$ws.Range("H133").Value = 56780
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 56780
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 56780
$ws.Range("N133").Value = -66900
$ws.Range("H138").Value = 1998.03
$ws.Range("I138").Value = 915.9706
$ws.Range("J138").Value = 2555.4546
$ws.Range("K138").Value = 2747.9118
$ws.Range("L138").Value = 7666.3638
$ws.Range("M138").Value = 2392.0882
$ws.Range("N138").Value = -17946.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2450.5625
$ws.Range("I2").Value = 2428.0908
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 2428.0908
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = -2315.0908
$ws.Range("N2").Value = -2726
$ws.Range("H45").Value = 2102.8
$ws.Range("I45").Value = 1900
$ws.Range("J45").Value = 2153.5
$ws.Range("K45").Value = 1900
$ws.Range("L45").Value = 2153.5
$ws.Range("M45").Value = -1523
$ws.Range("N45").Value = -2907.5
$ws.Range("H61").Value = 2122.1904
$ws.Range("I61").Value = 1944.2354
$ws.Range("J61").Value = 2878.5
$ws.Range("K61").Value = 1944.2354
$ws.Range("L61").Value = 2878.5
$ws.Range("M61").Value = -1732.2354
$ws.Range("N61").Value = -3302.5
$ws.Range("H63").Value = 3247.5
$ws.Range("I63").Value = 2000
$ws.Range("J63").Value = 4495
$ws.Range("K63").Value = 2000
$ws.Range("L63").Value = 4495
$ws.Range("M63").Value = -1314
$ws.Range("N63").Value = -5867
$ws.Range("H66").Value = 3247.5
$ws.Range("I66").Value = 2000
$ws.Range("J66").Value = 4495
$ws.Range("K66").Value = 10000
$ws.Range("L66").Value = 22475
$ws.Range("M66").Value = -6568
$ws.Range("N66").Value = -29339
$ws.Range("H116").Value = 2450.5625
$ws.Range("I116").Value = 2428.0908
$ws.Range("J116").Value = 2500
$ws.Range("K116").Value = 2428.0908
$ws.Range("L116").Value = 2500
$ws.Range("M116").Value = -134.0907999999999
$ws.Range("N116").Value = -7088
$ws.Range("H132").Value = 2701.4883
$ws.Range("I132").Value = 2486.6428
$ws.Range("J132").Value = 3102.5334
$ws.Range("K132").Value = 7459.928400000001
$ws.Range("L132").Value = 9307.600199999999
$ws.Range("M132").Value = -4929.928400000001
$ws.Range("N132").Value = -14367.6002
$ws.Range("H136").Value = 2122.1904
$ws.Range("I136").Value = 1944.2354
$ws.Range("J136").Value = 2878.5
$ws.Range("K136").Value = 5832.706200000001
$ws.Range("L136").Value = 8635.5
$ws.Range("M136").Value = -3282.706200000001
$ws.Range("N136").Value = -13735.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2450.5625
$ws.Range("I3").Value = 2428.0908
$ws.Range("J3").Value = 2500
$ws.Range("K3").Value = 2428.0908
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = -2314.0908
$ws.Range("N3").Value = -2728
$ws.Range("H86").Value = 2013.6364
$ws.Range("I86").Value = 1792.8572
$ws.Range("J86").Value = 2400
$ws.Range("K86").Value = 1792.8572
$ws.Range("L86").Value = 2400
$ws.Range("M86").Value = -669.8571999999999
$ws.Range("N86").Value = -4646
$ws.Range("H89").Value = 2013.6364
$ws.Range("I89").Value = 1792.8572
$ws.Range("J89").Value = 2400
$ws.Range("K89").Value = 8964.286
$ws.Range("L89").Value = 12000
$ws.Range("M89").Value = -3348.286
$ws.Range("N89").Value = -23232
$ws.Range("H92").Value = 110252
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 110252
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 110252
$ws.Range("N92").Value = -115244
$ws.Range("H93").Value = 30000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 30000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H94").Value = 1470.6923
$ws.Range("I94").Value = 1757
$ws.Range("J94").Value = 1136.6666
$ws.Range("K94").Value = 1757
$ws.Range("L94").Value = 1136.6666
$ws.Range("M94").Value = -1306
$ws.Range("N94").Value = -2038.6666
$ws.Range("H105").Value = 2387.9167
$ws.Range("I105").Value = 2321.6667
$ws.Range("J105").Value = 2454.1667
$ws.Range("K105").Value = 2321.6667
$ws.Range("L105").Value = 2454.1667
$ws.Range("M105").Value = -574.6667000000002
$ws.Range("N105").Value = -5948.1667
$ws.Range("H134").Value = 5793.8
$ws.Range("I134").Value = 5010.7
$ws.Range("J134").Value = 7360
$ws.Range("K134").Value = 15032.1
$ws.Range("L134").Value = 22080
$ws.Range("M134").Value = -12497.1
$ws.Range("N134").Value = -27150

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 810.84375
$ws.Range("I58").Value = 677.11536
$ws.Range("J58").Value = 1390.3334
$ws.Range("K58").Value = 677.11536
$ws.Range("L58").Value = 1390.3334
$ws.Range("M58").Value = -474.11536
$ws.Range("N58").Value = -1796.3334
$ws.Range("H132").Value = 3070.2856
$ws.Range("I132").Value = 2628.5881
$ws.Range("J132").Value = 4947.5
$ws.Range("K132").Value = 7885.7643
$ws.Range("L132").Value = 14842.5
$ws.Range("M132").Value = -5355.7643
$ws.Range("N132").Value = -19902.5
$ws.Range("H136").Value = 810.84375
$ws.Range("I136").Value = 677.11536
$ws.Range("J136").Value = 1390.3334
$ws.Range("K136").Value = 2031.34608
$ws.Range("L136").Value = 4171.0002
$ws.Range("M136").Value = 518.65392
$ws.Range("N136").Value = -9271.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 846.21
$ws.Range("I131").Value = 489
$ws.Range("J131").Value = 909.2471
$ws.Range("K131").Value = 1467
$ws.Range("L131").Value = 2727.7413
$ws.Range("M131").Value = 3573
$ws.Range("N131").Value = -12807.7413
$ws.Range("H132").Value = 3487914.8
$ws.Range("I132").Value = 1589167
$ws.Range("J132").Value = 12348738
$ws.Range("K132").Value = 14302503
$ws.Range("L132").Value = 111138642
$ws.Range("M132").Value = -14299973
$ws.Range("N132").Value = -111143702

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1344.8422
$ws.Range("I113").Value = 1074.2142
$ws.Range("J113").Value = 2102.6
$ws.Range("K113").Value = 1074.2142
$ws.Range("L113").Value = 2102.6
$ws.Range("M113").Value = 1095.7858
$ws.Range("N113").Value = -6442.6
$ws.Range("H132").Value = 3192.3125
$ws.Range("I132").Value = 3114.48
$ws.Range("J132").Value = 3470.2856
$ws.Range("K132").Value = 9343.440000000001
$ws.Range("L132").Value = 10410.8568
$ws.Range("M132").Value = -6813.440000000001
$ws.Range("N132").Value = -15470.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3418.44
$ws.Range("I7").Value = 3426.476
$ws.Range("J7").Value = 3376.25
$ws.Range("K7").Value = 3426.476
$ws.Range("L7").Value = 3376.25
$ws.Range("M7").Value = -3314.476
$ws.Range("N7").Value = -3600.25
$ws.Range("H26").Value = 9631.429
$ws.Range("I26").Value = 8350
$ws.Range("J26").Value = 11340
$ws.Range("K26").Value = 8350
$ws.Range("L26").Value = 11340
$ws.Range("M26").Value = -8055
$ws.Range("N26").Value = -11930
$ws.Range("H29").Value = 11471.286
$ws.Range("I29").Value = 6766.3335
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 6766.3335
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -6471.3335
$ws.Range("N29").Value = -15590
$ws.Range("H40").Value = 1921.3077
$ws.Range("I40").Value = 1329.6666
$ws.Range("J40").Value = 2428.4285
$ws.Range("K40").Value = 1329.6666
$ws.Range("L40").Value = 2428.4285
$ws.Range("M40").Value = -1193.6666
$ws.Range("N40").Value = -2700.4285
$ws.Range("H122").Value = 5138.095
$ws.Range("I122").Value = 5921.4287
$ws.Range("J122").Value = 3571.4285
$ws.Range("K122").Value = 17764.2861
$ws.Range("L122").Value = 10714.2855
$ws.Range("M122").Value = -15314.2861
$ws.Range("N122").Value = -15614.2855
$ws.Range("H126").Value = 3418.44
$ws.Range("I126").Value = 3426.476
$ws.Range("J126").Value = 3376.25
$ws.Range("K126").Value = 10279.428
$ws.Range("L126").Value = 10128.75
$ws.Range("M126").Value = -7809.428
$ws.Range("N126").Value = -15068.75
$ws.Range("H132").Value = 4644.364
$ws.Range("I132").Value = 5286.5
$ws.Range("J132").Value = 3873.8
$ws.Range("K132").Value = 15859.5
$ws.Range("L132").Value = 11621.4
$ws.Range("M132").Value = -13329.5
$ws.Range("N132").Value = -16681.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 936.25
$ws.Range("I81").Value = 858
$ws.Range("J81").Value = 1066.6666
$ws.Range("K81").Value = 1716
$ws.Range("L81").Value = 2133.3332
$ws.Range("M81").Value = -655
$ws.Range("N81").Value = -4255.3332
$ws.Range("H84").Value = 936.25
$ws.Range("I84").Value = 858
$ws.Range("J84").Value = 1066.6666
$ws.Range("K84").Value = 8580
$ws.Range("L84").Value = 10666.666
$ws.Range("M84").Value = -3276
$ws.Range("N84").Value = -21274.666
$ws.Range("H136").Value = 40003324
$ws.Range("I136").Value = 55559132
$ws.Range("J136").Value = 2672.8572
$ws.Range("K136").Value = 166677396
$ws.Range("L136").Value = 8018.571599999999
$ws.Range("M136").Value = -166674846
$ws.Range("N136").Value = -13118.5716
